$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header cells A1:C1 - the "16.5.1.1a." prefix lost its trailing period
#    (now a space) in front of the index name, in all three languages.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = '16.5.1.1a "Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк" индекси'
$ws.Range("B1").Value = '16.5.1.1a Индекс "Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления'''''
$ws.Range("C1").Value = '16.5.1.1a Index "Personal views about the level of corruption in executive government authorities and local government'''''

# ---------------------------------------------------------------------------
# 2) New column I - year 2020 data, mirroring the formatting of column H.
# ---------------------------------------------------------------------------

# Header (year) - same style as D4:H4
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020

# Top summary row (bold index row) - same style as H5
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 12.3
$ws.Range("I5").NumberFormat = "0.0"

# Region rows 6-13 - same style as H6:H13
$values = @{
    6  = 40.3
    7  = 36.2
    8  = 44.3
    9  = 36
    10 = 2.7
    11 = 32.9
    12 = 11.3
    13 = -18.2
}
foreach ($r in 6..13) {
    $ws.Range("H$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
    $ws.Range("I$r").Value = $values[$r]
}
$ws.Range("I6:I13").NumberFormat = "0.0"

# Bottom row (thick-bottom-bordered) - same style as H14
$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 33
$ws.Range("I14").NumberFormat = "0.0"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Selection moves from B16:B17 to the single cell F16.
# ---------------------------------------------------------------------------
$ws.Range("F16").Select() | Out-Null
